$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width (matches width="109.140625" in target)
$ws.Columns.Item(5).ColumnWidth = 109.140625

# Rows 61-71: Column F values
$fValues = @(2,3,4,5,6,7,8,11,12,1,9)
for ($i = 0; $i -lt $fValues.Length; $i++) {
    $row = 61 + $i
    $ws.Cells.Item($row, 6).Value2 = $fValues[$i]
}

# Rows 61-69: Column P values (same as F for those rows)
for ($i = 0; $i -lt 9; $i++) {
    $row = 61 + $i
    $ws.Cells.Item($row, 16).Value2 = $fValues[$i]
}

# Column E formulas rows 61-71
for ($row = 61; $row -le 71; $row++) {
    $ws.Cells.Item($row, 5).Formula = "=CONCATENATE(""insert into [UMCLOCKER].[dbo].[Locker](locker_index, locker_number,state,locker_type) values("",F$row,"",162,'AVAIABLE','F')"")"
}

# Q61 formula
$ws.Range("Q61").Formula = "=P61+P62"

# AutoFilter over A1:A79 -- apply before adding rows 80/81 so the persisted
# filter range stays fixed at A1:A79 instead of growing with new data rows
$ws.Range("A1:A79").AutoFilter()

# Defined name _xlnm._FilterDatabase, scoped to the sheet and hidden, same as
# what Excel creates automatically when an AutoFilter is applied
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$A`$79")
$fdb.Visible = $false

# New rows 80 and 81
$ws.Cells.Item(80, 1).Value2 = 162
$ws.Cells.Item(80, 2).Value2 = 162
$ws.Cells.Item(80, 3).Value2 = 12
$ws.Cells.Item(80, 4).Value2 = "Locker nữ"

$ws.Cells.Item(81, 1).Value2 = 165
$ws.Cells.Item(81, 2).Value2 = 165
$ws.Cells.Item(81, 3).Value2 = 12
$ws.Cells.Item(81, 4).Value2 = "Locker nữ"

# View settings: select E61:E71 and scroll so row 61 is the top-left
$ws.Range("E61:E71").Select()
$excel.ActiveWindow.ScrollRow = 61

$wb.Save()
